$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "Merging Microsoft Excel templates with OpenTBS has several big limitations because of the OpenXML format for Excel."
$ws.Range("B14").Value = "* Excel formulas are saved in Sheet subfiles with both the expression and the result. This means that when you write one TBS tag in Excel, it may have two in the XML source of the template."
$ws.Range("B15").Value = "* Texts in cells are not saved in the sheet subfile but in the 'xl/sharedStrings.xml' subfile. This means you cannot refer to any sheet tag (like row) for TBS tags placed in a text cell."
$ws.Range("B17").Value = "* Pictures placed in the sheet are not referenced in the sheet subfile but in another XML subfile (for instance 'xl\drawings\drawing1.xml'). This means you cannot use the usual parameter ""ope=changepic"" to change pictures in a sheet."
$ws.Range("B18").Value = "* Cells are saved in sheets with their absolute position (row+column). Thus, when you merge them using MergeBlock(), positions are duplicated. This can produce wrong sheets."
$ws.Range("B20").Value = "Far all those raisons, its seems than it is not possible to use MergeBlock() an Excel template."

$ws.Range("B22").Select()
